$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The author corrected the label in A7 from " 800 - 900" to " 800 -900"
# (removed the space before "900").
$ws.Range("A7").Value = " 800 -900"

# Reflect the last active selection being on A7 (as seen in the saved file).
$ws.Range("A7").Select()
